# Insert a new data row at row 335 (pushes the existing rows 335-366 down
# to 336-367, extending the used range to A1:R367) and populate it with a
# new weekly price record for "Haba" at Vega Central Mapocho de Santiago.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 335..366 down by one, creating a blank row 335.
$ws.Rows.Item(335).Insert()

# Populate the newly inserted row 335 with the new record.
$ws.Cells.Item(335, 1).Value = 9
$ws.Cells.Item(335, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(335, 3).Value = "Metropolitana"
$ws.Cells.Item(335, 4).Value = 45132
$ws.Cells.Item(335, 5).Value = 13
$ws.Cells.Item(335, 6).Value = 100112026
$ws.Cells.Item(335, 7).Value = "Haba"
$ws.Cells.Item(335, 8).Value = "Sin especificar"
$ws.Cells.Item(335, 9).Value = "Primera"
$ws.Cells.Item(335, 10).Value = 70
$ws.Cells.Item(335, 11).Value = 16000
$ws.Cells.Item(335, 12).Value = 17000
$ws.Cells.Item(335, 13).Value = 16500
$ws.Cells.Item(335, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(335, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(335, 16).Value = 660
$ws.Cells.Item(335, 17).Value = 25
$ws.Cells.Item(335, 18).Value = "Hortaliza"
